$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.544.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.90%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.589.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.62%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.60%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.26%  "

$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.593"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.606.44"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.11%  "

$ws.Range("E10").Value = "  +1.02%  "

$ws.Range("E11").Value = "  +1.61%  "

$ws.Range("E12").Value = "  +5.79%  "

$ws.Range("E13").Value = "  +4.94%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.050.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.56%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.515.35"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.70%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.04"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.22%  "

$ws.Range("E17").Value = "  +1.48%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.595.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.55%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.75%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "340.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.75%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.58"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.996"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.40%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.484"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.89%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "62.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.64%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.996"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.18%  "

$ws.Range("E27").Value = "  -1.19%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.49"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.89%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0770"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.08%  "

$ws.Range("E30").Value = "  -0.07%  "

$ws.Range("E31").Value = "  +0.59%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.12"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.59%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "157.53"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.91%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.40"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.84%  "

$ws.Range("E35").Value = "  +1.78%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.921"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.75%  "

$ws.Range("E37").Value = "  +3.28%  "

$ws.Range("E38").Value = "  +2.28%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.49"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.65%  "

$ws.Range("E40").Value = "  -4.10%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.68"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.10%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "290.02"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.60%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "136.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +9.21%  "

$ws.Range("E44").Value = "  -0.13%  "

$ws.Range("E45").Value = "  +0.41%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.598"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.99%  "

$ws.Range("E47").Value = "  +0.28%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0535"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.13%  "

$ws.Range("E49").Value = "  +1.91%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.977.06"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.92%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.65"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.50%  "
